$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of A2 so the new row matches existing rows
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Set the new value and move the active selection to the new cell
$ws.Range("A3").Value = "viji"
$ws.Range("A3").Select()
